# 2014-04-27 17:08
# "Conclusao das atividades da agenda" -- mark activities as concluded,
# unhide the previously-filtered rows, and clear the autofilter criteria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status updates for specific activities -----------------------------
# Row 32 (item 20): status goes from "?" (Em aberto) to "*" (Duvida)
$ws.Range("F32").Value = "*"

# Row 33 (item 21): collaborator changes from "Neimar" to "Aurelio"
$ws.Range("E33").Value = "Aurélio"

# Row 38 (item 26): status goes from "?" (Em aberto) to "!" (Concluida)
$ws.Range("F38").Value = "!"

# --- Unhide all the rows that were previously hidden by the autofilter --
for ($r = 4; $r -le 41; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- Clear the autofilter criteria on the "Colaborador" and "Status" ----
# columns (fields 5 and 6 of the A3:G43 filter range) so every row shows.
$ws.Range("A3:G43").AutoFilter(5)
$ws.Range("A3:G43").AutoFilter(6)

# --- Move the active selection to E3:F3 (anchored on F3) ----------------
$ws.Range("E3:F3").Select()
